$p = $ppt.ActivePresentation

# --- Slide 2: merge the three runs of the "House renting..." paragraph into one ---
$s2 = $p.Slides.Item(2)
$sh2 = $s2.Shapes.Item("Content Placeholder 2")
$tr2 = $sh2.TextFrame.TextRange
$para1 = $tr2.Paragraphs(1, 1)

$run1Len = "House renting is a web application that is people ".Length
$tail = $tr2.Characters($para1.Start + $run1Len, $para1.Length - $run1Len)
$tail.Delete() | Out-Null

$run1 = $tr2.Characters($para1.Start, $run1Len)
$run1.InsertAfter("can advertise their own house and people can view house advertisement who want to rent house.") | Out-Null

# --- Slide 6: split the "Actually, this project is popular in urban area." run into three ---
$s6 = $p.Slides.Item(6)
$sh6 = $s6.Shapes.Item("Content Placeholder 2")
$tr6 = $sh6.TextFrame.TextRange
$para = $tr6.Paragraphs(2, 1)
$start = $para.Start
$len = $para.Length

$whole = $tr6.Characters($start, $len)
$whole.Text = "Actually, this project will be popular in urban area."

$para2 = $tr6.Paragraphs(2, 1)
$t1 = "Actually, this "
$t2 = "project will be "
$t3 = "popular in urban area."

$c1 = $tr6.Characters($para2.Start, $t1.Length)
$c1.Text = $t1

$c2 = $tr6.Characters($para2.Start + $t1.Length, $t2.Length)
$c2.Text = $t2

$c3 = $tr6.Characters($para2.Start + $t1.Length + $t2.Length, $t3.Length)
$c3.Text = $t3
